# Daily attendance processing - 2025-11-21 02:59:07
#
# For every data row in the "Recorded By" column, rotate the comma-separated
# list of recorders so the last entry moves to the front (e.g.
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"). Rows whose
# list already starts with "System" (i.e. already rotated / single entry)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Locate the "Recorded By" column from the header row; fall back to column G (7).
$targetCol = 7
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Value()
    if ($header -ne $null -and $header.Equals("Recorded By")) {
        $targetCol = $c
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $targetCol)
    $val = $cell.Value()
    if ($val -eq $null) { continue }

    $parts = $val -split ", "
    if ($parts.Length -le 1) { continue }

    $firstElem = $parts[0]
    if ($firstElem.Equals("System")) { continue }

    $lastElem = $parts[$parts.Length - 1]
    $rest = $parts[0..($parts.Length - 2)]
    $newParts = @($lastElem) + $rest
    $newVal = $newParts -join ", "

    $cell.Value = $newVal
}
